$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns F, G, H (row 1), matching the style of the
# existing header cells (e.g. C1) so they look the same (bold, bordered,
# centered). Copy the format from an existing header cell first, then set
# the text, so the new cells reuse the same cell style as A1:E1.
$ws.Range("C1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Boolean outlier flags for rows 2-11: all FALSE
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 6).Value = $false
    $ws.Cells.Item($r, 7).Value = $false
    $ws.Cells.Item($r, 8).Value = $false
}

# Row 12: F12 = TRUE, G12 = FALSE, H12 = FALSE
$ws.Cells.Item(12, 6).Value = $true
$ws.Cells.Item(12, 7).Value = $false
$ws.Cells.Item(12, 8).Value = $false
